# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" positioned between "总计" and "2021-Q3".
# - Populate "2022-Q4" with the fund holding table (fund 002174).
# - Update the "总计" sheet: row 2 now reflects the newest quarter (2022-Q4),
#   and a new row 3 preserves the previous quarter's totals (2021-Q3).

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)   # "总计"
$q3Sheet    = $wb.Worksheets.Item(2)   # "2021-Q3" (currently 2nd, used as the
                                        # insertion anchor for the new sheet)

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right before the existing "2021-Q3" sheet
#    so the tab order becomes: 总计, 2022-Q4, 2021-Q3
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Fill in the "2022-Q4" sheet contents
# ---------------------------------------------------------------------------
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("C2").Value = "东方互联网嘉混合"
$q4Sheet.Range("H2").Value = 8

# B2, D2, E2, F2, G2 hold numeric-looking text (leading zero / fixed
# decimals) that must stay text, e.g. "002174" not 2174. A direct .Value
# assignment auto-coerces these into numbers. To keep them as genuine text
# WITHOUT minting a new cell style (NumberFormat="@" would create one),
# stage them as formula results on a scratch row, then paste only the
# VALUES into place - pasted values keep the source's text-ness but none of
# its formatting, so no new style is added to the workbook.
$scratch = $q4Sheet.Range("B100:G100")
$q4Sheet.Range("B100").Formula = '="002174"'
$q4Sheet.Range("C100").Formula = '="东方互联网嘉混合"'
$q4Sheet.Range("D100").Formula = '="1.17"'
$q4Sheet.Range("E100").Formula = '="93.41"'
$q4Sheet.Range("F100").Formula = '="1.44"'
$q4Sheet.Range("G100").Formula = '="0.0168"'
$scratch.Copy()
$q4Sheet.Range("B2:G2").PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()

# Formatting: the header row + A2 use the same style already present on the
# "总计" sheet's header row (bold + thin border + centered), so copy it
# instead of re-building it (re-building would mint a brand-new style).
$totalSheet.Range("B1:D1").Copy()
$q4Sheet.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats
$totalSheet.Range("B1:D1").Copy()
$q4Sheet.Range("E1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Update the "总计" sheet: row 2 -> 2022-Q4 figures, new row 3 -> the old
#    2021-Q3 figures that used to live in row 2.
# ---------------------------------------------------------------------------
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = 0.02

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.01

# A3 carries the same style as A2 (bold/border/center).
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) The "2021-Q3" sheet (now 3rd) was the selected tab before this edit and
#    its own content/selection is untouched by this change, so restore it as
#    the active tab (adding the new sheet made it active by default).
#    Re-fetch the sheet reference by name since it shifted position.
# ---------------------------------------------------------------------------
$q3SheetNow = $wb.Worksheets.Item("2021-Q3")
$q3SheetNow.Select() | Out-Null
$q3SheetNow.Range("A1").Select() | Out-Null
